$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DepartmentData")

$ws.Range("A12").Value = "ciphggd Dep"
$ws.Range("B12").Value = "DEP0071"

$ws.Range("A13").Value = "lkdpjno Dep"
$ws.Range("B13").Value = "DEP0073"

$ws.Range("A14").Value = "fjbkfdg Dep"
$ws.Range("B14").Value = "DEP0074"
